# Update the date title
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-09 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-09-10 Wednesday", 2)

# Update the division-problem table. Addressed by (row, column) so the
# two cells that happen to share the same original text ("34÷5=6, 4")
# each get their own, distinct replacement.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "95÷9=10, 5" },
    @{ Row = 1;  Col = 2; Text = "56÷3=18, 2" },
    @{ Row = 1;  Col = 3; Text = "43÷9=4, 7" },
    @{ Row = 1;  Col = 4; Text = "74÷8=9, 2" },
    @{ Row = 1;  Col = 5; Text = "41÷8=5, 1" },

    @{ Row = 5;  Col = 1; Text = "50÷5=10, 0" },
    @{ Row = 5;  Col = 2; Text = "84÷3=28, 0" },
    @{ Row = 5;  Col = 3; Text = "25÷7=3, 4" },
    @{ Row = 5;  Col = 4; Text = "87÷4=21, 3" },
    @{ Row = 5;  Col = 5; Text = "36÷5=7, 1" },

    @{ Row = 9;  Col = 1; Text = "28÷4=7, 0" },
    @{ Row = 9;  Col = 2; Text = "82÷6=13, 4" },
    @{ Row = 9;  Col = 3; Text = "70÷4=17, 2" },
    @{ Row = 9;  Col = 4; Text = "62÷3=20, 2" },
    @{ Row = 9;  Col = 5; Text = "36÷7=5, 1" },

    @{ Row = 13; Col = 1; Text = "57÷5=11, 2" },
    @{ Row = 13; Col = 2; Text = "88÷4=22, 0" },
    @{ Row = 13; Col = 3; Text = "74÷2=37, 0" },
    @{ Row = 13; Col = 4; Text = "68÷7=9, 5" },
    @{ Row = 13; Col = 5; Text = "90÷2=45, 0" },

    @{ Row = 17; Col = 1; Text = "26÷6=4, 2" },
    @{ Row = 17; Col = 2; Text = "28÷6=4, 4" },
    @{ Row = 17; Col = 3; Text = "96÷3=32, 0" },
    @{ Row = 17; Col = 4; Text = "44÷9=4, 8" },
    @{ Row = 17; Col = 5; Text = "41÷4=10, 1" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
